$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 43-49: narrative notes introducing the FWHM analysis, column A ---
$ws.Range("A43").Value = 'The analysis results shown below are for looking for the average FWHM for the scan data of the structures. I use the same Jupyter Notebook as above.'
$ws.Range("A44").Value = 'scipy.signal.peak_widths is used to do this with the parameter rel_height=0.5. This means peak width is found at h=h_peak-PR where h_peak is the height of the peaks (note the peaks are found by flipping the resonance data because the resonances are minima).'
$ws.Range("A45").Value = 'scipy.signal.peak_widths is found in samples, and multiplied by wavelength_step_size to give the FWHM in nm'
$ws.Range("A46").Value = 'peak data is found in same way as above'
$ws.Range("A47").Value = 'definitions unless stated otherwise are same as above. I will also add the fsr data below as I might as well, because notebook calculates both fsr and FWHM data.'
$ws.Range("A48").Value = 'The parameter wlen passed into peak_widths is approx_fsr/wavelength step size'
$ws.Range("A49").Value = 'I will Git commit after each data run and entry into this table.'

# --- Row 53: header row for the new FWHM results table (mirrors row 26 + 2 extra cols) ---
$ws.Range("A53").Value = 'Data CSV Filename'
$ws.Range("B53").Value = 'Wavelength step size/nm'
$ws.Range("C53").Value = 'Start array index'
$ws.Range("D53").Value = 'End array index'
$ws.Range("E53").Value = 'Start wavelength/nm'
$ws.Range("F53").Value = 'End wavelength/nm'
$ws.Range("G53").Value = 'prominence/dBm'
$ws.Range("H53").Value = 'distance'
$ws.Range("I53").Value = 'approx_fsr/nm'
$ws.Range("J53").Value = 'fsr_mean/nm'
$ws.Range("K53").Value = 'fsr_std error/nm'
$ws.Range("L53").Value = 'double count check passed?'
$ws.Range("M53").Value = 'mean FWHM/nm'
$ws.Range("N53").Value = 'FWHM error/nm'

# --- Column widths: widen columns whose header/content got longer (approximation of Excel AutoFit) ---
$ws.Columns.Item(5).ColumnWidth = 17.666666
$ws.Columns.Item(6).ColumnWidth = 16.833333
$ws.Columns.Item(7).ColumnWidth = 14.666666
$ws.Columns.Item(11).ColumnWidth = 14
$ws.Columns.Item(13).ColumnWidth = 14.5
$ws.Columns.Item(14).ColumnWidth = 14

# --- Leave the view scrolled near the new table, matching where the author was working ---
$ws.Range("A48").Select()
